$d = $word.ActiveDocument

function Get-ParagraphByText($searchText) {
    # Search the whole document body for the first paragraph whose text
    # matches $searchText exactly, and return that paragraph object.
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $searchText"
    }
    return $rng.Paragraphs(1)
}

# ---------------------------------------------------------------------
# The "COMPETENCES TECHNIQUES" skill list is being reordered:
#
#   Before:                                          After:
#   Langages : ...                                   Langages : ...
#   Bases de données : SQL, MongoDB, Neo4j, Redis    Visualisation : tableau
#   Autres : tensorflow dataiku, databricks          MLOps : Git, DVC, ...
#   Visualisation : tableau                          Autres : tensorflow dataiku, databricks
#   ML/AI : Scikit-Learn, ...                        ML/AI : Scikit-Learn, ...
#   MLOps : Git, DVC, ...                            Bases de données : SQL, MongoDB, Neo4j, Redis
#
# i.e. "Visualisation" and "MLOps" move up right after "Langages", and
# "Bases de données" moves down to become the last line. "Autres" and
# "ML/AI" keep their relative position.
#
# We apply this using only text replacements plus one paragraph
# insertion/deletion pair, walking the affected paragraphs via
# Next()/Previous() navigation (rather than re-running Find on text that
# becomes temporarily duplicated mid-edit) so each step unambiguously
# targets the correct paragraph.
# ---------------------------------------------------------------------

# Step 1: the paragraph that used to read "Bases de données : ..." becomes
# "Visualisation : tableau", and a brand-new paragraph carrying the
# "MLOps : ..." text is inserted right after it (pushing "Autres" etc. down).
$pBases = Get-ParagraphByText("Bases de données : SQL, MongoDB, Neo4j, Redis")
$pBases.Range.Text = "Visualisation : tableau"
$pBases.Range.InsertParagraphAfter()
$pNewMLOps = $pBases.Next()
$pNewMLOps.Range.Text = "MLOps : Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit"

# Step 2: "Autres : ..." is unaffected; move to it, then to the old
# "Visualisation : tableau" paragraph right after it, and turn that one
# into the "ML/AI : ..." text.
$pAutres = $pNewMLOps.Next()
$pOldVisualisation = $pAutres.Next()
$pOldVisualisation.Range.Text = "ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn"

# Step 3: the following paragraph (the old "ML/AI : ...") becomes the new
# last "Bases de données : ..." line.
$pOldMLAI = $pOldVisualisation.Next()
$pOldMLAI.Range.Text = "Bases de données : SQL, MongoDB, Neo4j, Redis"

# Step 4: the trailing paragraph after that (the original "MLOps : ...",
# now duplicated by Step 1's insertion) is removed entirely.
$pOldMLOps = $pOldMLAI.Next()
$pOldMLOps.Range.Delete()
